$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet/tab from "General_Surgery" to "Session"
$ws.Name = "Session"

# Remove the duplicate log entry that was re-scanned/re-entered at row 20
# (Student ID 191258 logged at 11:27:36). Deleting the entire row shifts
# every subsequent row up by one, which also updates the used range from
# A1:F111 down to A1:F110 and keeps all remaining rows/columns intact.
$ws.Rows(20).Delete()
